# Apply edit: insert a new data row at row 4 (shifting existing rows 4-46 down to 5-47),
# and populate the new row 4 with its data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts rows 4..46 down to 5..47
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new record's data
$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(4, 3).Value = "Maule"
$ws.Cells.Item(4, 4).Value = 44496
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 6).Value = 100112026
$ws.Cells.Item(4, 7).Value = "Haba"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 500
$ws.Cells.Item(4, 11).Value = 7000
$ws.Cells.Item(4, 12).Value = 7000
$ws.Cells.Item(4, 13).Value = 7000
$ws.Cells.Item(4, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región del Maule"
$ws.Cells.Item(4, 16).Value = 280
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format style used by the other date cells in column D
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
